# Fruta / hortaliza, semanal
# Shuffle the weekly "Fecha" / "Volumen" / "Precio" block (columns D, M, N, O, P, S)
# across rows 2-20 of the data table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current D/M/N/O/P/S values for every data row before
# writing anything, since several rows exchange values with each other.
$cols = @("D", "M", "N", "O", "P", "S")
$snapshot = @{}
for ($r = 2; $r -le 20; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $row
}

# Destination row -> source row (the row whose D/M/N/O/P/S block moves in).
$mapping = @{
    2  = 5
    3  = 13
    4  = 14
    5  = 3
    6  = 2
    7  = 18
    8  = 19
    9  = 15
    10 = 8
    11 = 20
    12 = 12
    13 = 9
    14 = 6
    15 = 17
    16 = 10
    17 = 11
    18 = 16
    19 = 7
    20 = 4
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcData[$c]
    }
}
